# Update the "EZ Water Adjustment" sheet's water-chemistry inputs and
# recipe-selection state to match the new revision of the recipe.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EZ Water Adjustment")

# Starting water profile: Chloride went from 0 to 71 ppm.
$ws.Range("G5").Value = 71

# Mash water additions row (grams / ml):
#  - Calc. Chloride source ratio recomputed (2.45 -> 2.38 g/ml)
#  - Epsom Salt acid content 0.86 -> 0.83
#  - Lactic Acid ml 0 -> 1
$ws.Range("E37").Formula = "=2.38/1.335"
$ws.Range("F37").Value = 0.83
$ws.Range("J37").Value = 1

# Restore the sheet's active selection to a single cell (J38), undoing the
# prior multi-cell selection/scroll position.
$ws.Activate()
$ws.Range("J38").Select()
